# Update countries & provincias Spain
# Rebuild the full country data table (A4:H216) with the corrected
# country ordering (sorted descending by "Casos totales") and refreshed
# case counts, per the updated daily snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @"
Estados Unidos|644417|328|48708|567150|13487|30|28559
España|182816|2157|74797|88889|7371|318|19130
Italia|165155|0|38092|105418|3079|0|21645
Francia|147863|0|30955|99741|6457|0|17167
Alemania|134753|0|77000|53949|4288|0|3804
Reino Unido|98476|0|0|85264|1559|0|12868
China|82341|46|77892|1107|95|0|3342
Iran|77995|1606|52229|20897|3594|92|4869
Turquia|69392|0|5674|62200|1820|0|1518
Belgica|34809|1236|7562|22390|1182|417|4857
Paises Bajos|29214|1061|250|25649|1279|181|3315
Brasil|29015|405|14026|13229|296|3|1760
Canada|28379|0|8979|18390|557|0|1010
Rusia|27938|3448|2304|25402|8|34|232
Suiza|26422|86|15400|9753|386|30|1269
Portugal|18841|750|493|17719|229|30|629
Austria|14412|62|8098|5921|232|0|393
Israel|12591|90|2624|9827|174|10|140
Irlanda|12547|0|77|12026|158|0|444
Suecia|12540|613|381|10826|996|130|1333
India|12456|86|1513|10520|0|1|423
Peru|11475|0|3108|8113|146|0|254
Corea del Sur|10613|22|7757|2627|55|4|229
Japon|8626|0|901|7547|168|0|178
Chile|8273|0|2937|5242|389|0|94
Ecuador|7858|0|780|6690|135|0|388
Polonia|7771|189|774|6705|160|6|292
Rumania|7707|491|1357|5958|243|20|392
Pakistan|6919|536|1645|5146|46|17|128
Dinamarca|6879|198|3023|3535|89|12|321
Noruega|6798|1|32|6616|64|0|150
Australia|6468|21|3747|2658|66|0|63
Arabia Saudita|6380|518|990|5307|71|4|83
Chequia|6303|2|831|5306|75|0|166
Mexico|5847|448|2125|3273|207|43|449
Filipinas|5660|207|435|4863|1|13|362
Indonesia|5516|380|548|4472|0|27|496
Emiratos Arabes Unidos|5365|0|1034|4298|1|0|33
Malasia|5182|110|2766|2332|56|1|84
Serbia|4873|0|400|4374|128|0|99
Bielorrusia|4204|476|203|3961|65|4|40
Ucrania|4161|397|186|3859|45|8|116
Catar|4103|392|415|3681|37|0|7
Panama|3751|0|75|3573|106|0|103
Singapur|3699|0|652|3037|29|0|10
Republica Dominicana|3614|0|208|3217|143|0|189
Luxemburgo|3373|0|526|2778|33|0|69
Finlandia|3369|132|300|2997|75|0|72
Colombia|3105|0|452|2522|106|0|131
Tailandia|2672|29|1593|1033|61|3|46
Argentina|2571|0|631|1825|117|3|115
Sudafrica|2506|0|410|2062|7|0|34
Egipto|2505|0|589|1733|0|0|183
Marruecos|2251|227|247|1876|1|1|128
Grecia|2192|0|269|1821|72|0|102
Argelia|2160|0|708|1116|60|0|336
Moldavia|2049|0|235|1763|80|5|51
Croacia|1791|50|529|1227|31|2|35
Islandia|1727|0|1077|642|8|0|8
Barein|1698|27|703|988|3|0|7
Hungria|1652|73|199|1311|60|8|142
Banglades|1572|341|49|1463|1|10|60
Kuwait|1524|119|225|1296|32|0|3
Estonia|1434|34|133|1265|10|1|36
Irak|1415|0|812|524|0|0|79
Nueva Zelanda|1401|15|770|622|3|0|9
Kazajistan|1362|67|273|1072|22|1|17
Uzbekistan|1349|47|129|1216|8|0|4
Eslovenia|1268|20|174|1033|31|0|61
Azerbaiyan|1253|0|404|836|24|0|13
Bosnia y Herzegovina|1167|57|277|847|4|2|43
Armenia|1159|48|358|783|30|1|18
Lituania|1128|37|178|920|14|0|30
Oman|1019|109|176|839|3|0|4
Hong Kong|1018|1|485|529|9|0|4
Eslovaquia|977|114|167|802|5|2|8
Republica de Macedonia|974|0|98|831|15|0|45
Camerun|848|0|165|666|0|0|17
Afganistan|840|56|54|756|0|5|30
Cuba|814|0|151|639|15|0|24
Bulgaria|783|36|122|624|31|1|37
Tunez|780|0|43|702|89|0|35
Republica de Chipre|715|0|65|638|8|0|12
Crucero|712|0|644|56|7|0|12
Letonia|675|9|57|613|3|0|5
Principado de Andorra|673|0|169|471|17|0|33
Libano|663|5|85|557|30|0|21
Costa de Marfil|654|0|146|502|0|0|6
Ghana|641|0|83|550|2|0|8
Costa Rica|626|0|67|555|11|0|4
Niger|584|0|90|480|0|0|14
Burkina Faso|542|0|226|284|0|0|32
Albania|518|24|277|215|5|1|26
Uruguay|493|0|272|212|13|0|9
Kirguistan|466|17|91|370|5|0|5
Bolivia|441|44|14|398|3|1|29
Republica de Yibuti|435|0|71|362|0|0|2
Honduras|426|7|9|382|10|4|35
Malta|412|13|82|327|4|0|3
Nigeria|407|0|128|267|2|0|12
Guinea|404|0|31|372|0|0|1
Jordania|401|0|250|144|5|0|7
Taiwan|395|0|155|234|0|0|6
San Marino|393|0|53|304|15|0|36
Reunion|391|0|237|154|3|0|0
Estado de Palestina|374|0|63|309|0|0|2
Georgia|336|30|74|259|6|0|3
Mauricio|324|0|65|250|3|0|9
Senegal|314|0|190|122|1|0|2
Montenegro|290|2|55|231|7|0|4
Isla de Man|283|27|153|126|13|0|4
Vietnam|268|0|176|92|8|0|0
Consejo Danes para los Refugiados|267|13|23|222|0|1|22
Sri Lanka|238|0|68|163|1|0|7
Mayotte|233|16|69|161|3|0|3
Kenia|225|0|53|162|2|0|10
Venezuela|197|0|111|77|6|0|9
Guatemala|196|16|19|172|3|0|5
Islas Feroe|184|0|169|15|0|0|0
Paraguay|174|13|30|136|1|0|8
El Salvador|164|5|33|125|2|0|6
Martinica|158|0|73|77|17|0|8
Mali|148|0|34|101|0|0|13
Guadalupe|145|0|67|70|13|0|8
Ruanda|136|0|54|82|0|0|0
Brunei|136|0|108|27|2|0|1
Gibraltar|131|0|104|27|1|0|0
Jamaica|125|0|21|99|0|0|5
Camboya|122|0|98|24|1|0|0
Congo|117|0|11|101|0|0|5
Trinidad yTobago|114|0|20|86|0|0|8
Madagascar|111|1|33|78|1|0|0
Tanzania|94|6|11|79|0|0|4
Monaco|93|0|12|78|2|0|3
Aruba|93|0|39|53|1|0|1
Etiopia|92|7|15|74|0|0|3
Guayana Francesa|86|0|51|35|1|0|0
Birmania|85|11|2|79|0|0|4
Togo|81|0|35|43|0|0|3
Bermudas|81|0|33|43|3|0|5
Gabon|80|0|4|75|0|0|1
Somalia|80|0|2|73|2|0|5
Liechtenstein|79|0|55|23|0|0|1
Barbados|75|2|15|55|4|0|5
Islas Caimanes|60|0|6|53|3|0|1
Liberia|59|0|4|49|0|0|6
San Martin (Parte Holandesa)|57|4|12|36|6|0|9
Cabo Verde|56|0|1|54|0|0|1
Polinesia Francesa|55|0|0|55|1|0|0
Uganda|55|0|12|43|0|0|0
Guyana|55|0|8|41|5|0|6
Bahamas|53|0|6|39|1|0|8
Guinea Ecuatorial|51|0|4|47|0|0|0
Libia|48|0|11|36|0|0|1
Zambia|48|0|30|16|1|0|2
Macao|45|0|16|29|1|0|0
Guinea-Bisau|43|0|0|43|0|0|0
Haiti|41|0|0|38|0|0|3
Puerto Rico|39|0|1|36|0|0|2
Eritrea|35|0|0|35|0|0|0
San Martin (Parte Francesa)|35|0|13|20|5|0|2
Benin|35|0|18|16|0|0|1
Siria|33|0|5|26|0|0|2
Guam|32|0|0|31|0|0|1
Sudan|32|0|4|23|0|0|5
Mongolia|31|1|5|26|0|0|0
Mozambique|29|0|2|27|0|0|0
Republica del Chad|27|4|5|22|0|0|0
Zimbabue|23|0|1|19|0|0|3
Antigua y Barbuda|23|0|3|17|1|1|3
Maldivas|23|1|16|7|0|0|0
Laos|19|0|2|17|0|0|0
Angola|19|0|5|12|0|0|2
Timor Oriental|18|10|1|17|0|0|0
Nueva Caledonia|18|0|1|17|1|0|0
Belice|18|0|0|16|1|0|2
Fiyi|17|1|0|17|0|0|0
Islas Virgenes de los Estados Unidos|17|0|0|17|0|0|0
Nepal|16|0|1|15|0|0|0
Malaui|16|0|0|14|1|0|2
Namibia|16|0|3|13|0|0|0
Dominica|16|0|8|8|0|0|0
Suazilandia|16|1|8|8|0|0|0
Sierra Leona|15|2|0|15|0|0|0
Botsuana|15|2|0|14|0|0|1
Santa Lucia|15|0|11|4|0|0|0
Granada|14|0|0|14|2|0|0
San Cristobal y Nieves|14|0|0|14|0|0|0
Curazao|14|0|10|3|0|0|1
San Vicente y las Granadinas|12|0|1|11|0|0|0
Republica de Africa Central|12|0|4|8|0|0|0
Seychelles|11|0|0|11|0|0|0
Islas Malvinas|11|0|1|10|0|0|0
Montserrat|11|0|1|10|1|0|0
Groenlandia|11|0|11|0|0|0|0
Islas Turcas y Caicos|10|0|0|9|0|0|1
Surinam|10|0|6|3|0|0|1
Gambia|9|0|2|6|0|0|1
Nicaragua|9|0|4|4|0|0|1
Santa Sede|8|0|2|6|0|0|0
Mauritania|7|0|2|4|0|0|1
Sahara Occidental|6|0|0|6|0|0|0
San Bartolome|6|0|4|2|0|0|0
Burundi|5|0|0|4|0|0|1
Butan|5|0|2|3|0|0|0
Santo Tome y Principe|4|0|0|4|0|0|0
Sudan del Sur|4|0|0|4|0|0|0
Bonaire, San Eustaquio y Saba|3|0|0|3|0|0|0
Anguila|3|0|1|2|0|0|0
Islas Virgenes Britanicas|3|0|2|1|0|0|0
Papua Nueva Guinea|2|0|0|2|0|0|0
San Pedro y Miquelon|1|0|0|1|0|0|0
Yemen|1|0|0|1|0|0|0
"@

$lines = $data -split "`n"
$rows = New-Object System.Collections.ArrayList
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    [void]$rows.Add(($line -split "\|"))
}

$n = $rows.Count
$arr = New-Object 'object[,]' $n,8

for ($i = 0; $i -lt $n; $i++) {
    $fields = $rows[$i]
    $arr[$i,0] = $fields[0]
    for ($j = 1; $j -lt 8; $j++) {
        $arr[$i,$j] = [double]$fields[$j]
    }
}

$startRow = 4
$endRow = $startRow + $n - 1
$rng = $ws.Range("A" + $startRow + ":H" + $endRow)
$rng.Value = $arr
